$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3304176666666667
$ws.Range("H2").Value = 0.9912529999999999
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3412459999999999
$ws.Range("N2").Value = 1.023738
$ws.Range("O2").Value = 0.147839631507836
$ws.Range("P2").Value = 0.147839631507836
$ws.Range("Q2").Value = 0.1127537070793333
$ws.Range("R2").Value = 1.014783363714
$ws.Range("S2").Value = 0.147839631507836
$ws.Range("T2").Value = 0.147839631507836
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3304176666666667
$ws.Range("H3").Value = 0.9912529999999999
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4973860000000001
$ws.Range("N3").Value = 1.492158
$ws.Range("O3").Value = 0.2154849081224587
$ws.Range("P3").Value = 0.2154849081224587
$ws.Range("Q3").Value = 0.1643451215526667
$ws.Range("R3").Value = 1.479106093974
$ws.Range("S3").Value = 0.2154849081224587
$ws.Range("T3").Value = 0.2154849081224587
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3304176666666667
$ws.Range("H4").Value = 0.9912529999999999
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.299649
$ws.Range("N4").Value = 0.8989469999999999
$ws.Range("O4").Value = 0.1298183648795636
$ws.Range("P4").Value = 0.1298183648795636
$ws.Range("Q4").Value = 0.099009323399
$ws.Range("R4").Value = 0.8910839105909999
$ws.Range("S4").Value = 0.1298183648795636
$ws.Range("T4").Value = 0.1298183648795636
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3304176666666667
$ws.Range("H5").Value = 0.9912529999999999
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4296976666666667
$ws.Range("N5").Value = 1.289093
$ws.Range("O5").Value = 0.186159968760885
$ws.Range("P5").Value = 0.186159968760885
$ws.Range("Q5").Value = 0.1419797003921111
$ws.Range("R5").Value = 1.277817303529
$ws.Range("S5").Value = 0.186159968760885
$ws.Range("T5").Value = 0.186159968760885
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3304176666666667
$ws.Range("H6").Value = 0.9912529999999999
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1086676666666667
$ws.Range("N6").Value = 0.326003
$ws.Range("O6").Value = 0.04707861131505237
$ws.Range("P6").Value = 0.04707861131505237
$ws.Range("Q6").Value = 0.03590571686211111
$ws.Range("R6").Value = 0.323151451759
$ws.Range("S6").Value = 0.04707861131505237
$ws.Range("T6").Value = 0.04707861131505237
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3304176666666667
$ws.Range("H7").Value = 0.9912529999999999
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.631571
$ws.Range("N7").Value = 1.894713
$ws.Range("O7").Value = 0.2736185154142042
$ws.Range("P7").Value = 0.2736185154142042
$ws.Range("Q7").Value = 0.2086822161543333
$ws.Range("R7").Value = 1.878139945389
$ws.Range("S7").Value = 0.2736185154142042
$ws.Range("T7").Value = 0.2736185154142042
